$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.921.61"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "'1.876.70"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'0.7430"
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("D6").Value = "'242.59"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "'0.07249"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'24.69"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").Value = "'0.08406"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "'0.7524"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "'5.437"
$ws.Range("D14").Value = "'1.894.04"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").Value = "'92.61"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "'29.919.17"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'6.086"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "'247.96"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "'13.59"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'0.000007864"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'2.125.52"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "'8.039"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("D26").Value = "'9.273"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'165.21"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "'18.65"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "'2.039"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'1.510"
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("D31").Value = "'4.612"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "'1.533"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "'4.281"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "'0.05335"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'0.7538"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").Value = "'0.9982"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'2.692"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'0.01965"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "'0.4531"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").Value = "'1.112.23"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'6.051"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'72.71"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "'0.8545"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "'103.44"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "'1.860"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "'7.631"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.023.43"
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.906"
$ws.Range("E51").Value = "  -2.62%  "
